# Update cryptos list values (price + 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.943.25"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "'2.446.89"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'569.79"
$ws.Range("D6").Value = "'146.74"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.534"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").Value = "'2.449.65"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "'26.90"
$ws.Range("E14").Value = "  +3.06%  "
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").Value = "'2.901.31"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").Value = "'62.881.40"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "'2.450.08"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "'11.38"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "'7.24"
$ws.Range("E20").Value = "  +6.38%  "
$ws.Range("D21").Value = "'324.47"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").Value = "  +12.13%  "
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").Value = "'619.02"
$ws.Range("E26").Value = "  +10.44%  "
$ws.Range("D27").Value = "'8.59"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("D28").Value = "'0.0000102"
$ws.Range("E28").Value = "  +9.95%  "
$ws.Range("D29").Value = "'2.568.61"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'1.47"
$ws.Range("E31").Value = "  +6.07%  "
$ws.Range("D32").Value = "'8.28"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").Value = "'1.90"
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("D35").Value = "'5.08"
$ws.Range("E35").Value = "  +7.54%  "
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D39").Value = "'5.39"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("D40").Value = "'18.68"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").Value = "'144.65"
$ws.Range("E41").Value = "  -5.04%  "
$ws.Range("D42").Value = "'1.78"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("E43").Value = "  +16.36%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'146.99"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").Value = "'0.0539"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").Value = "'20.61"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").Value = "'0.601"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").Value = "'0.0921"
$ws.Range("E51").Value = "  -0.04%  "
